# Refresh market-price / leve-profit figures (currentAveragePrice*, Leve*Price, Leve*Profit)
# for a handful of leves across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets, as produced by the
# scheduled market-data runner. Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ.
# M/N are only present when the corresponding NQ/HQ leve price (K/L) is nonzero, so some
# writes clear a cell instead of setting it to 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -877
$ws.Range("N86").ClearContents() | Out-Null

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4384
$ws.Range("N89").ClearContents() | Out-Null

# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 2489.8
$ws.Range("I92").Value = 2500.4285
$ws.Range("K92").Value = 2500.4285
$ws.Range("M92").Value = -1252.4285

# Row 96: Scroll Down / Grade 1 Reisui of Intelligence
$ws.Range("H96").Value = 1878.5
$ws.Range("I96").Value = 838
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 2514
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -1141
$ws.Range("N96").Value = -17746

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 16668346
$ws.Range("I100").Value = 18183558
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 18183558
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = -18183017
$ws.Range("N100").Value = -2082

# Row 123: Nearly Bare / Gaja Grimoire
$ws.Range("H123").Value = 41890
$ws.Range("J123").Value = 41890
$ws.Range("L123").Value = 41890
$ws.Range("N123").Value = -51690

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2283.7144
$ws.Range("I137").Value = 1558.6
$ws.Range("J137").Value = 4096.5
$ws.Range("K137").Value = 4675.799999999999
$ws.Range("L137").Value = 12289.5
$ws.Range("M137").Value = -2125.799999999999
$ws.Range("N137").Value = -17389.5

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3295.04
$ws.Range("I138").Value = 640.96875
$ws.Range("J138").Value = 4544.0146
$ws.Range("K138").Value = 1922.90625
$ws.Range("L138").Value = 13632.0438
$ws.Range("M138").Value = 3217.09375
$ws.Range("N138").Value = -23912.0438

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1143.25
$ws.Range("I45").Value = 1179.4166
$ws.Range("J45").Value = 1034.75
$ws.Range("K45").Value = 1179.4166
$ws.Range("L45").Value = 1034.75
$ws.Range("M45").Value = -802.4166
$ws.Range("N45").Value = -1788.75

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1011.86365
$ws.Range("I110").Value = 1064.4
$ws.Range("J110").Value = 486.5
$ws.Range("K110").Value = 1064.4
$ws.Range("L110").Value = 486.5
$ws.Range("M110").Value = 980.5999999999999
$ws.Range("N110").Value = -4576.5

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 3809.8125
$ws.Range("I99").Value = 1451.6
$ws.Range("J99").Value = 4881.727
$ws.Range("K99").Value = 1451.6
$ws.Range("L99").Value = 4881.727
$ws.Range("M99").Value = 46.40000000000009
$ws.Range("N99").Value = -7877.727

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1557.1111
$ws.Range("I105").Value = 1554
$ws.Range("J105").Value = 1596
$ws.Range("K105").Value = 1554
$ws.Range("L105").Value = 1596
$ws.Range("M105").Value = 193
$ws.Range("N105").Value = -5090

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 7408668
$ws.Range("I16").Value = 13889965
$ws.Range("K16").Value = 13889965
$ws.Range("M16").Value = -13889678

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 14708381
$ws.Range("I31").Value = 1251.7727
$ws.Range("J31").Value = 41671450
$ws.Range("K31").Value = 1251.7727
$ws.Range("L31").Value = 41671450
$ws.Range("M31").Value = -956.7727
$ws.Range("N31").Value = -41672040

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 14708381
$ws.Range("I34").Value = 1251.7727
$ws.Range("J34").Value = 41671450
$ws.Range("K34").Value = 1251.7727
$ws.Range("L34").Value = 41671450
$ws.Range("M34").Value = -1049.7727
$ws.Range("N34").Value = -41671854

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 6974
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6974
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6974
$ws.Range("M62").ClearContents() | Out-Null
$ws.Range("N62").Value = -8222

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 6974
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6974
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 34870
$ws.Range("M65").ClearContents() | Out-Null
$ws.Range("N65").Value = -41110

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 2318.7273
$ws.Range("I105").Value = 1919.5
$ws.Range("J105").Value = 2797.8
$ws.Range("K105").Value = 1919.5
$ws.Range("L105").Value = 2797.8
$ws.Range("M105").Value = -172.5
$ws.Range("N105").Value = -6291.8

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 7408668
$ws.Range("I113").Value = 13889965
$ws.Range("K113").Value = 13889965
$ws.Range("M113").Value = -13887795

$ws = $wb.Worksheets.Item("CUL")
# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 2347.0483
$ws.Range("I121").Value = 199
$ws.Range("J121").Value = 2418.65
$ws.Range("K121").Value = 597
$ws.Range("L121").Value = 7255.950000000001
$ws.Range("M121").Value = 713
$ws.Range("N121").Value = -9875.950000000001

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 9616180
$ws.Range("I131").Value = 62500230
$ws.Range("J131").Value = 897.0682
$ws.Range("K131").Value = 187500690
$ws.Range("L131").Value = 2691.2046
$ws.Range("M131").Value = -187495650
$ws.Range("N131").Value = -12771.2046

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1144.2106
$ws.Range("I61").Value = 1043.5294
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1043.5294
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -841.5293999999999
$ws.Range("N61").Value = -2404

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1144.2106
$ws.Range("I113").Value = 1043.5294
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1043.5294
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1126.4706
$ws.Range("N113").Value = -6340

# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 33850
$ws.Range("J133").Value = 33850
$ws.Range("L133").Value = 33850
$ws.Range("N133").Value = -38910

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 290.5
$ws.Range("J100").Value = 350
$ws.Range("L100").Value = 700
$ws.Range("N100").Value = -1782

# Row 125: Color Coated / Almasty Serge Coat of Healing
$ws.Range("H125").Value = 49143.332
$ws.Range("J125").Value = 49143.332
$ws.Range("L125").Value = 49143.332
$ws.Range("N125").Value = -58983.332

# Row 128: Lightening Up / Scarlet Moko Gaskins of the Rising Dragon
$ws.Range("H128").Value = 41855
$ws.Range("J128").Value = 41855
$ws.Range("L128").Value = 41855
$ws.Range("N128").Value = -51815

# Row 131: A Better Bottom Line / AR-Caean Velvet Bottoms of Scouting
$ws.Range("H131").Value = 67903.336
$ws.Range("J131").Value = 67903.336
$ws.Range("L131").Value = 67903.336
$ws.Range("N131").Value = -77983.336

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 4918.636
$ws.Range("I136").Value = 2600.8572
$ws.Range("J136").Value = 8974.75
$ws.Range("K136").Value = 7802.571599999999
$ws.Range("L136").Value = 26924.25
$ws.Range("M136").Value = -5252.571599999999
$ws.Range("N136").Value = -32024.25
